# Restructure the "Customer" currency lookup sheet:
#  - drop the old helper column C (raw currency names) entirely
#  - column A now holds the currency label directly (no more
#    SUBSTITUTE(..) formula rebuilding the old "Values" text)
#  - column B's formula gets its result passed through an extra
#    SUBSTITUTE so the generated code token has no inner spaces
#    (e.g. "US Dollar" -> "USD" instead of "US_D")
#  - the old "Values" header is no longer needed, only "Label" /
#    "BOM to XOM" remain

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("A1").Value = "Label"
$ws.Range("B1").Value = "BOM to XOM"

# Column A: plain currency names (used to be rebuilt via
# SUBSTITUTE(C#, " ", "_") from column C -- now column C is gone and
# the label itself lives directly in A)
$labels = @("US Dollar", "Euro", "Pound", "JP Yen", "PRC RMB", "SG Dollar")
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}

# Column B: same CONCATENATE idea as before, but the 4-letter prefix
# now also has its spaces stripped out via an inner SUBSTITUTE
$ws.Range("B2").Formula = '=CONCATENATE("return """, SUBSTITUTE(UPPER(LEFT(A2,4)), " ", ""), """;")'
$ws.Range("B3:B7").Formula = '=CONCATENATE("return """, SUBSTITUTE(UPPER(LEFT(A3,4)), " ", ""), """;")'

# Column C (old raw "Values" helper column) is no longer needed
$ws.Columns("C").Delete()

# Leave the selection where the author left it after deleting the
# now-empty column C
$ws.Columns("C").Select()
